# CellAlignment demo: rename sheet, update label text, clear stray sample
# cells, resize the row, and left/center-align the remaining cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "Sheet0"

# Clear the extra sample cells (A3:C3 numbers, E3 boolean) leaving only D3.
$ws.Range("A3:C3").ClearContents()
$ws.Range("E3").ClearContents()

# Update the surviving string cell's text and alignment.
$ws.Range("D3").Value = "Align It"
$ws.Range("D3").HorizontalAlignment = -4131
$ws.Range("D3").VerticalAlignment = -4108

# Make row 3 taller to show off the vertical centering.
$ws.Rows.Item(3).RowHeight = 50
